$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (prices like "63.298.51") are stored as text, not auto-converted to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.298.51'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.649.64'
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '607.71'
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('D6').Value = '144.04'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '2.648.79'
$ws.Range('E9').Value = '  +3.63%  '
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('D11').Value = '5.64'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '0.364'
$ws.Range('E13').Value = '  +4.81%  '
$ws.Range('D14').Value = '27.37'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '3.122.96'
$ws.Range('E15').Value = '  +3.87%  '
$ws.Range('D16').Value = '63.151.30'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '2.659.96'
$ws.Range('E18').Value = '  +4.52%  '
$ws.Range('D19').Value = '11.42'
$ws.Range('E19').Value = '  +3.60%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '342.89'
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = '4.44'
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('E22').Value = '  +4.19%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '67.12'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +2.80%  '
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('D27').Value = '8.66'
$ws.Range('E27').Value = '  +6.30%  '
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').Value = '548.66'
$ws.Range('E29').Value = '  +17.67%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').Value = '2.04'
$ws.Range('E32').Value = '  +5.16%  '
$ws.Range('E33').Value = '  +8.46%  '
$ws.Range('D34').Value = '0.0₃0811'
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').Value = '172.29'
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('D36').Value = '5.11'
$ws.Range('E36').Value = '  +13.59%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  +2.50%  '
$ws.Range('D39').Value = '19.14'
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('E40').Value = '  +8.50%  '
$ws.Range('D41').Value = '171.97'
$ws.Range('E41').Value = '  +9.14%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = '3.77'
$ws.Range('E43').Value = '  +2.26%  '
$ws.Range('D44').Value = '22.37'
$ws.Range('E44').Value = '  +5.15%  '
$ws.Range('D45').Value = '0.0578'
$ws.Range('E45').Value = '  +8.52%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '18.85'
$ws.Range('E49').Value = '  +5.22%  '
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('D51').Value = '11.22'
$ws.Range('E51').Value = '  -1.21%  '

# Restore default (unstyled) cell style now that the text values are set,
# matching the original workbook formatting.
$ws.Range("B2:E51").Style = "Normal"
